$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 180
$ws.Range("I2").Value = 486
$ws.Range("J2").Value = 2091
$ws.Range("K2").Value = 15
$ws.Range("L2").Value = 550
$ws.Range("M2").Value = 39
$ws.Range("N2").Value = 344
$ws.Range("P2").Value = 5
$ws.Range("R2").Value = 43
$ws.Range("S2").Value = 245
$ws.Range("T2").Value = 364
$ws.Range("V2").Value = 3257
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 3222
$ws.Range("Y2").Value = 3
$ws.Range("Z2").Value = 49
$ws.Range("AA2").Value = 18
